# Apply cell value updates to Sheet1 per the Flashscore odds update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("I6").Value = 2.9
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7
$ws.Range("AY6").Value = 29
# Row 7
$ws.Range("W7").Value = 4.75
$ws.Range("AE7").Value = 26
$ws.Range("AH7").Value = 19
$ws.Range("AI7").Value = 17
$ws.Range("AQ7").Value = 51
# Row 8
$ws.Range("M8").Value = 1.03
$ws.Range("P8").Value = 3.5
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.93
# Row 9
$ws.Range("G9").Value = 2.38
$ws.Range("I9").Value = 2.7
$ws.Range("M9").Value = 1.03
$ws.Range("O9").Value = 1.19
$ws.Range("AA9").Value = 19
$ws.Range("AN9").Value = 4.75
$ws.Range("AZ9").Value = 41
# Row 13
$ws.Range("N13").Value = 15
# Row 14
$ws.Range("N14").Value = 8
$ws.Range("W14").Value = 6
# Row 16
$ws.Range("G16").Value = 2.35
$ws.Range("H16").Value = 3.1
$ws.Range("I16").Value = 2.88
$ws.Range("K16").Value = 2.05
$ws.Range("L16").Value = 3.6
$ws.Range("N16").Value = 8.5
$ws.Range("Q16").Value = 2.15
$ws.Range("R16").Value = 1.67
$ws.Range("U16").Value = 1.83
$ws.Range("V16").Value = 1.83
$ws.Range("W16").Value = 7.5
$ws.Range("X16").Value = 11
$ws.Range("AB16").Value = 34
$ws.Range("AC16").Value = 8.5
$ws.Range("AD16").Value = 6
$ws.Range("AH16").Value = 15
$ws.Range("AK16").Value = 26
$ws.Range("AM16").Value = 301
$ws.Range("AN16").Value = 4.33
$ws.Range("AP16").Value = 26
$ws.Range("AS16").Value = 201
$ws.Range("AY16").Value = 29
# Row 17
$ws.Range("G17").Value = 2.63
$ws.Range("I17").Value = 2.5
$ws.Range("J17").Value = 3.4
$ws.Range("L17").Value = 3.25
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 9
$ws.Range("W17").Value = 8
$ws.Range("X17").Value = 13
$ws.Range("Y17").Value = 11
$ws.Range("AA17").Value = 23
$ws.Range("AC17").Value = 9
$ws.Range("AH17").Value = 12
$ws.Range("AI17").Value = 10
$ws.Range("AJ17").Value = 26
$ws.Range("AK17").Value = 21
$ws.Range("AN17").Value = 4.75
$ws.Range("AR17").Value = 81
$ws.Range("AW17").Value = 4.5
$ws.Range("BA17").Value = 67
# Row 22
$ws.Range("Q22").Value = 2.15
$ws.Range("R22").Value = 1.67
# Row 27
$ws.Range("H27").Value = 3.15
$ws.Range("I27").Value = 3.05
$ws.Range("L27").Value = 3.6
$ws.Range("T27").Value = 2.72
$ws.Range("W27").Value = 7.6
$ws.Range("Y27").Value = 8.75
$ws.Range("AA27").Value = 18
$ws.Range("AB27").Value = 28
$ws.Range("AG27").Value = 9.25
$ws.Range("AK27").Value = 27
$ws.Range("AL27").Value = 35
$ws.Range("AN27").Value = 4.15
$ws.Range("AT27").Value = 2.72
$ws.Range("AX27").Value = 17
$ws.Range("AY27").Value = 24
$ws.Range("AZ27").Value = 80
# Row 35
$ws.Range("I35").Value = 4.25
$ws.Range("L35").Value = 4.5
$ws.Range("P35").Value = 2.9
$ws.Range("Q35").Value = 1.95
$ws.Range("V35").Value = 1.82
$ws.Range("AG35").Value = 11.25
$ws.Range("AH35").Value = 24
$ws.Range("AL35").Value = 45
$ws.Range("AM35").Value = 600
$ws.Range("AV35").Value = 60
$ws.Range("AY35").Value = 28
$ws.Range("AZ35").Value = 120
$ws.Range("BA35").Value = 150
# Row 36
$ws.Range("G36").Value = 2.95
$ws.Range("I36").Value = 2.35
$ws.Range("J36").Value = 3.4
$ws.Range("L36").Value = 3
$ws.Range("Q36").Value = 1.88
$ws.Range("S36").Value = 1.39
$ws.Range("T36").Value = 2.55
$ws.Range("X36").Value = 16.5
$ws.Range("Y36").Value = 10.25
$ws.Range("Z36").Value = 37
$ws.Range("AA36").Value = 24
$ws.Range("AB36").Value = 28
$ws.Range("AD36").Value = 6.1
$ws.Range("AE36").Value = 12.5
$ws.Range("AK36").Value = 19.5
$ws.Range("AL36").Value = 27
$ws.Range("AN36").Value = 4.85
$ws.Range("AO36").Value = 15.5
$ws.Range("AP36").Value = 21
$ws.Range("AQ36").Value = 70
$ws.Range("AR36").Value = 90
$ws.Range("AW36").Value = 4.3
$ws.Range("AZ36").Value = 55
$ws.Range("BA36").Value = 90
